$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = "11-50" + [char]10 + "(up to 100 for construction sector)"
$ws.Range("B21").Value = "51-250" + [char]10 + "(up to 400 for construction sector)"
$ws.Range("B22").Value = ">250" + [char]10 + "(> 400 for construction sector)"
